$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal string (e.g. "0.9993") must be
# forced to stay text (matching the original inline-string cell type) --
# otherwise Excel auto-converts them to numbers on assignment. We flip the
# cell to a Text number format just for the write, then clear the format
# again so the cell keeps its original (default) style.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "24.358.38"
$ws.Range("E2").Value = "  -2.23%  "

$ws.Range("D3").Value = "1.647.26"
$ws.Range("E3").Value = "  -4.67%  "

$ws.Range("E4").Value = "  -1.42%  "

Set-TextValue $ws.Range("D5") "0.9993"
$ws.Range("E5").Value = "  -1.00%  "

Set-TextValue $ws.Range("D6") "305.71"
$ws.Range("E6").Value = "  -2.58%  "

Set-TextValue $ws.Range("D7") "0.3628"
$ws.Range("E7").Value = "  -4.28%  "

Set-TextValue $ws.Range("D8") "47.37"
$ws.Range("E8").Value = "  -4.49%  "

Set-TextValue $ws.Range("D9") "0.3267"
$ws.Range("E9").Value = "  -7.56%  "

Set-TextValue $ws.Range("D10") "1.113"
$ws.Range("E10").Value = "  -7.35%  "

Set-TextValue $ws.Range("D11") "0.06900"
$ws.Range("E11").Value = "  -8.50%  "

Set-TextValue $ws.Range("D12") "1.002"
$ws.Range("E12").Value = "  -1.42%  "

Set-TextValue $ws.Range("D13") "5.943"
$ws.Range("E13").Value = "  -7.19%  "

Set-TextValue $ws.Range("D14") "19.09"
$ws.Range("E14").Value = "  -8.95%  "

$ws.Range("D15").Value = "1.645.40"
$ws.Range("E15").Value = "  -5.19%  "

Set-TextValue $ws.Range("D16") "6.545"
$ws.Range("E16").Value = "  -6.86%  "

$ws.Range("E17").Value = "  -8.10%  "

Set-TextValue $ws.Range("D18") "0.06486"
$ws.Range("E18").Value = "  -3.50%  "

Set-TextValue $ws.Range("D19") "1.000"
$ws.Range("E19").Value = "  -1.01%  "

Set-TextValue $ws.Range("D20") "76.71"
$ws.Range("E20").Value = "  -9.99%  "

Set-TextValue $ws.Range("D21") "5.887"
$ws.Range("E21").Value = "  -8.49%  "

Set-TextValue $ws.Range("D22") "15.67"
$ws.Range("E22").Value = "  -9.73%  "

Set-TextValue $ws.Range("D23") "12.14"
$ws.Range("E23").Value = "  -6.58%  "

$ws.Range("D24").Value = "24.367.41"
$ws.Range("E24").Value = "  -2.61%  "

Set-TextValue $ws.Range("D25") "2.413"
$ws.Range("E25").Value = "  -1.64%  "

Set-TextValue $ws.Range("D26") "2.333"
$ws.Range("E26").Value = "  -17.27%  "

Set-TextValue $ws.Range("D27") "145.51"
$ws.Range("E27").Value = "  -4.60%  "

Set-TextValue $ws.Range("D28") "18.42"
$ws.Range("E28").Value = "  -10.73%  "

$ws.Range("D29").Value = "1.828.26"
$ws.Range("E29").Value = "  -5.64%  "

Set-TextValue $ws.Range("D30") "124.58"
$ws.Range("E30").Value = "  -5.88%  "

Set-TextValue $ws.Range("D31") "1.151"
$ws.Range("E31").Value = "  -2.85%  "

Set-TextValue $ws.Range("D32") "4.045"
$ws.Range("E32").Value = "  -4.61%  "

Set-TextValue $ws.Range("D33") "5.553"
$ws.Range("E33").Value = "  -20.33%  "

Set-TextValue $ws.Range("D34") "0.08311"
$ws.Range("E34").Value = "  -5.23%  "

Set-TextValue $ws.Range("D35") "1.668"
$ws.Range("E35").Value = "  -7.48%  "

Set-TextValue $ws.Range("D36") "12.24"
$ws.Range("E36").Value = "  -11.73%  "

Set-TextValue $ws.Range("D37") "5.175"
$ws.Range("E37").Value = "  -8.34%  "

Set-TextValue $ws.Range("D38") "0.06025"
$ws.Range("E38").Value = "  -9.03%  "

Set-TextValue $ws.Range("D39") "0.02204"
$ws.Range("E39").Value = "  -10.17%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D40") "1.203"
$ws.Range("E40").Value = "  -5.65%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "8.210"
$ws.Range("E41").Value = "  -11.07%  "

Set-TextValue $ws.Range("D42") "0.2033"
$ws.Range("E42").Value = "  -8.24%  "

Set-TextValue $ws.Range("D43") "0.9997"
$ws.Range("E43").Value = "  -1.05%  "

$ws.Range("E44").Value = "  -10.07%  "

Set-TextValue $ws.Range("D45") "3.709"
$ws.Range("E45").Value = "  -4.28%  "

Set-TextValue $ws.Range("D46") "12.56"
$ws.Range("E46").Value = "  -11.02%  "

Set-TextValue $ws.Range("D47") "0.5567"
$ws.Range("E47").Value = "  -10.09%  "

Set-TextValue $ws.Range("D48") "121.51"
$ws.Range("E48").Value = "  -6.57%  "

$ws.Range("E49").Value = "  -10.38%  "

Set-TextValue $ws.Range("D50") "0.06878"
$ws.Range("E50").Value = "  -5.58%  "

Set-TextValue $ws.Range("D51") "73.72"
$ws.Range("E51").Value = "  -7.67%  "
